# Selenium testdata workbook: add login test cases
# (verifyThatAdminCanLogInWithValidCredentials / verifyThatAdminCannotLogInWithInvalidCredentials)
# replacing the old loginLogoutTest / newTest rows.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------------
# RUNMANAGER sheet: rename the two test cases (descriptions/flags untouched)
# ---------------------------------------------------------------------------
$wsRun = $wb.Worksheets.Item("RUNMANAGER")
$wsRun.Range("A2").Value = "verifyThatAdminCanLogInWithValidCredentials"
$wsRun.Range("A3").Value = "verifyThatAdminCannotLogInWithInvalidCredentials"

$wsRun.Columns.Item(1).ColumnWidth = 41
$wsRun.Range("E3").Select() | Out-Null

# ---------------------------------------------------------------------------
# DATA sheet: rebuild the 2 test x 2 browser login matrix
# ---------------------------------------------------------------------------
$wsData = $wb.Worksheets.Item("DATA")

# Drop the old 6th row (loginLogoutTest/subscribe); the new matrix only
# needs 4 data rows.
$wsData.Rows.Item(6).Delete()

$wsData.Range("A2").Value = "verifyThatAdminCanLogInWithValidCredentials"
$wsData.Range("B2").Value = "yes"
$wsData.Range("C2").Value = "Admin"
$wsData.Range("D2").Value = "admin123"
$wsData.Range("E2").Value = "amuthan"
$wsData.Range("F2").Value = "chrome"

$wsData.Range("A3").Value = "verifyThatAdminCanLogInWithValidCredentials"
$wsData.Range("B3").Value = "yes"
$wsData.Range("C3").Value = "Admin"
$wsData.Range("D3").Value = "admin123"
$wsData.Range("E3").Value = "amuthan"
$wsData.Range("F3").Value = "firefox"

$wsData.Range("A4").Value = "verifyThatAdminCannotLogInWithInvalidCredentials"
$wsData.Range("B4").Value = "yes"
$wsData.Range("C4").Value = "admin12"
$wsData.Range("D4").Value = "admin123"
$wsData.Range("E4").Value = "sunil"
$wsData.Range("F4").Value = "chrome"

$wsData.Range("A5").Value = "verifyThatAdminCannotLogInWithInvalidCredentials"
$wsData.Range("B5").Value = "yes"
$wsData.Range("C5").Value = "admin12"
$wsData.Range("D5").Value = "admin123"
$wsData.Range("E5").Value = "sunil"
$wsData.Range("F5").Value = "firefox"

$wsData.Columns.Item(1).ColumnWidth = 41
$wsData.Range("F5").Select() | Out-Null

# RUNMANAGER stays the active/visible tab, as in the source workbook.
$wsRun.Activate()
